$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$m = Get-Member -InputObject $ws.ListObjects
Write-Host ($m | Out-String)
